# Update "想去人数" (interest count) and "最低票价" (lowest price) figures
# across all four sheets to match the refreshed GitHub Pages data export.
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 41757
$ws.Range("G2").Value = 85
$ws.Range("F5").Value = 9467
$ws.Range("F7").Value = 848
$ws.Range("G7").Value = 61.2
$ws.Range("F8").Value = 891
$ws.Range("F9").Value = 721
$ws.Range("F10").Value = 211
$ws.Range("F12").Value = 293
$ws.Range("F13").Value = 890
$ws.Range("F15").Value = 123
$ws.Range("F16").Value = 728
$ws.Range("F17").Value = 313
$ws.Range("F18").Value = 1394
$ws.Range("F20").Value = 651
$ws.Range("F21").Value = 696
$ws.Range("F22").Value = 453
$ws.Range("F23").Value = 682
$ws.Range("F24").Value = 725
$ws.Range("F27").Value = 61
$ws.Range("F28").Value = 496
$ws.Range("F29").Value = 519
$ws.Range("F30").Value = 49
$ws.Range("F31").Value = 237
$ws.Range("F32").Value = 924
$ws.Range("F35").Value = 91
$ws.Range("F36").Value = 211
$ws.Range("F38").Value = 385
$ws.Range("F39").Value = 1251
$ws.Range("F40").Value = 289
$ws.Range("F42").Value = 1227
$ws.Range("F43").Value = 373
$ws.Range("F46").Value = 30

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 213
$ws.Range("F5").Value = 4444
$ws.Range("F7").Value = 330
$ws.Range("F10").Value = 77
$ws.Range("F19").Value = 4385

# ---- Sheet 3: 本地生活 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2026
$ws.Range("F3").Value = 517
$ws.Range("F4").Value = 390

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2026
$ws.Range("F3").Value = 517
$ws.Range("F4").Value = 41757
$ws.Range("G4").Value = 85
$ws.Range("F7").Value = 213
$ws.Range("F8").Value = 330
$ws.Range("F10").Value = 9467
$ws.Range("F12").Value = 848
$ws.Range("G12").Value = 61.2
$ws.Range("F13").Value = 848
$ws.Range("G13").Value = 61.2
$ws.Range("F14").Value = 77
$ws.Range("F15").Value = 390
$ws.Range("F16").Value = 891
$ws.Range("F18").Value = 211
$ws.Range("F19").Value = 293
$ws.Range("F20").Value = 890
$ws.Range("F24").Value = 728
$ws.Range("F25").Value = 313
$ws.Range("F26").Value = 1394
$ws.Range("F27").Value = 651
$ws.Range("F28").Value = 696
$ws.Range("F29").Value = 453
$ws.Range("F30").Value = 682
$ws.Range("F31").Value = 725
$ws.Range("F33").Value = 61
$ws.Range("F34").Value = 496
$ws.Range("F35").Value = 49
$ws.Range("F36").Value = 237
$ws.Range("F37").Value = 924
$ws.Range("F40").Value = 92
$ws.Range("F41").Value = 211
$ws.Range("F42").Value = 385
$ws.Range("F43").Value = 1227
$ws.Range("F44").Value = 373
$ws.Range("F47").Value = 30
